$d = $word.ActiveDocument

# The title paragraph is the first paragraph in the document (style "Title"),
# e.g. "Visualizing the ocean floor". We need to turn the trailing "floor"
# run into "floor:" and then append a dynamic-title suffix of
# " word_document" where "word_document" carries the VerbatimChar style,
# e.g. "Visualizing the ocean floor: word_document".
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Style.NameLocal -eq "Title") {
        $titlePara = $candidate
        break
    }
}
if ($titlePara -eq $null) {
    $titlePara = $d.Paragraphs(1)
}
$titleRange = $titlePara.Range

# Find the last run's extent ("floor") without touching any of the
# surrounding runs -- locate it via Find restricted to the title's range so
# occurrences of the same word elsewhere in the document are untouched.
$findRange = $d.Range($titleRange.Start, $titleRange.End)
$findRange.Find.ClearFormatting()
$findRange.Find.Execute("floor", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

$wordStart = $findRange.Start
$wordEnd = $findRange.End

# Replace just the matched word's text with an empty string, then type the
# new text back in via InsertAfter from a collapsed range. This keeps the
# run separate from its unchanged neighbours (a plain Range.Text assignment
# would otherwise coalesce it with the preceding run).
$wordRange = $d.Range($wordStart, $wordEnd)
$wordRange.Text = ""

$newWordRun = $d.Range($wordStart, $wordStart)
$newWordRun.InsertAfter("floor:")

# Append a single space as its own run right after "floor:", still inside
# the title paragraph (just before the paragraph mark).
$afterColon = $titlePara.Range.End - 1
$spaceInsert = $d.Range($afterColon, $afterColon)
$spaceInsert.InsertAfter(" ")

# Append "word_document" as its own run, then style it with the Verbatim
# character style so it renders like inline code.
$afterSpace = $titlePara.Range.End - 1
$verbStart = $afterSpace
$verbInsert = $d.Range($afterSpace, $afterSpace)
$verbInsert.InsertAfter("word_document")
$verbEnd = $titlePara.Range.End - 1

$verbRange = $d.Range($verbStart, $verbEnd)
$verbRange.Style = "VerbatimChar"

Write-Output $titlePara.Range.Text
